# Scheduled-runner refresh of market price / profit figures across the
# Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# for the rows whose underlying market data changed; clears cells that
# no longer have data and adds cells that newly do.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 826.5
$ws.Range("I125").Value = 606.38464
$ws.Range("J125").Value = 1235.2858
$ws.Range("K125").Value = 5457.46176
$ws.Range("L125").Value = 11117.5722
$ws.Range("M125").Value = -2997.46176
$ws.Range("N125").Value = -16037.5722

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 3000
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H32").Value = 1521778.2
$ws.Range("I32").Value = 1547737.1
$ws.Range("J32").Value = 120000
$ws.Range("K32").Value = 1547737.1
$ws.Range("L32").Value = 120000
$ws.Range("M32").Value = -1547450.1
$ws.Range("N32").Value = -120574

$ws.Range("H97").Value = 721.0333000000001
$ws.Range("I97").Value = 686.8214
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 686.8214
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -190.8214
$ws.Range("N97").Value = -2192

$ws.Range("H102").Value = 8580
$ws.Range("I102").Value = 2230
$ws.Range("K102").Value = 2230
$ws.Range("M102").Value = -608

$ws.Range("H122").Value = 1773
$ws.Range("I122").Value = 538
$ws.Range("J122").Value = 5889.6665
$ws.Range("K122").Value = 1614
$ws.Range("L122").Value = 17668.9995
$ws.Range("M122").Value = 836
$ws.Range("N122").Value = -22568.9995

$ws.Range("H132").Value = 35794.2
$ws.Range("I132").Value = 60519.824
$ws.Range("J132").Value = 3460.6924
$ws.Range("K132").Value = 181559.472
$ws.Range("L132").Value = 10382.0772
$ws.Range("M132").Value = -179029.472
$ws.Range("N132").Value = -15442.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6753.846
$ws.Range("I86").Value = 7066.6665
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 7066.6665
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -5943.6665
$ws.Range("N86").Value = -5246

$ws.Range("H89").Value = 6753.846
$ws.Range("I89").Value = 7066.6665
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 35333.3325
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -29717.3325
$ws.Range("N89").Value = -26232

$ws.Range("H94").Value = 918.125
$ws.Range("I94").Value = 821.53845
$ws.Range("J94").Value = 1336.6666
$ws.Range("K94").Value = 821.53845
$ws.Range("L94").Value = 1336.6666
$ws.Range("M94").Value = -370.53845
$ws.Range("N94").Value = -2238.6666

$ws.Range("H105").Value = 49221.45
$ws.Range("I105").Value = 87724.55
$ws.Range("J105").Value = 2162.111
$ws.Range("K105").Value = 87724.55
$ws.Range("L105").Value = 2162.111
$ws.Range("M105").Value = -85977.55
$ws.Range("N105").Value = -5656.111

$ws.Range("H107").Value = 1628.3334
$ws.Range("I107").Value = 1773.8182
$ws.Range("J107").Value = 1228.25
$ws.Range("K107").Value = 1773.8182
$ws.Range("L107").Value = 1228.25
$ws.Range("M107").Value = 146.1818000000001
$ws.Range("N107").Value = -5068.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1059514.2
$ws.Range("I31").Value = 836.35895
$ws.Range("J31").Value = 1937991.6
$ws.Range("K31").Value = 836.35895
$ws.Range("L31").Value = 1937991.6
$ws.Range("M31").Value = -541.35895
$ws.Range("N31").Value = -1938581.6

$ws.Range("H34").Value = 1059514.2
$ws.Range("I34").Value = 836.35895
$ws.Range("J34").Value = 1937991.6
$ws.Range("K34").Value = 836.35895
$ws.Range("L34").Value = 1937991.6
$ws.Range("M34").Value = -634.35895
$ws.Range("N34").Value = -1938395.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 99
$ws.Range("I34").Value = 99
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 297
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -213
$ws.Range("N34").ClearContents()

$ws.Range("H39").Value = 5975
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 5975
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 17925
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -18513

$ws.Range("H55").Value = 3089.5
$ws.Range("I55").Value = 400
$ws.Range("J55").Value = 3388.3333
$ws.Range("K55").Value = 1200
$ws.Range("L55").Value = 10164.9999
$ws.Range("M55").Value = -1023
$ws.Range("N55").Value = -10518.9999

$ws.Range("H68").Value = 1083.6897
$ws.Range("I68").Value = 570.9706
$ws.Range("J68").Value = 1810.0416
$ws.Range("K68").Value = 1712.9118
$ws.Range("L68").Value = 5430.1248
$ws.Range("M68").Value = -901.9117999999999
$ws.Range("N68").Value = -7052.1248

$ws.Range("H71").Value = 1083.6897
$ws.Range("I71").Value = 570.9706
$ws.Range("J71").Value = 1810.0416
$ws.Range("K71").Value = 5138.7354
$ws.Range("L71").Value = 16290.3744
$ws.Range("M71").Value = -1082.7354
$ws.Range("N71").Value = -24402.3744

$ws.Range("H92").Value = 487.75
$ws.Range("I92").Value = 317.33334
$ws.Range("K92").Value = 952.0000200000001
$ws.Range("M92").Value = 295.9999799999999

$ws.Range("H107").Value = 853
$ws.Range("I107").Value = 282.4138
$ws.Range("J107").Value = 1723.8948
$ws.Range("K107").Value = 847.2413999999999
$ws.Range("L107").Value = 5171.6844
$ws.Range("M107").Value = 1072.7586
$ws.Range("N107").Value = -9011.6844

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4858
$ws.Range("I102").Value = 3185.7778
$ws.Range("J102").Value = 6739.25
$ws.Range("K102").Value = 3185.7778
$ws.Range("L102").Value = 6739.25
$ws.Range("M102").Value = -1563.7778
$ws.Range("N102").Value = -9983.25

$ws.Range("H107").Value = 5892
$ws.Range("I107").Value = 7891.231
$ws.Range("J107").Value = 694
$ws.Range("K107").Value = 7891.231
$ws.Range("L107").Value = 694
$ws.Range("M107").Value = -5971.231
$ws.Range("N107").Value = -4534

$ws.Range("H140").Value = 38677.25
$ws.Range("I140").Value = 17709
$ws.Range("J140").Value = 45666.668
$ws.Range("K140").Value = 17709
$ws.Range("L140").Value = 45666.668
$ws.Range("M140").Value = -12529
$ws.Range("N140").Value = -56026.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 8833.333000000001
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 8833.333000000001
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 8833.333000000001
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -9519.333000000001

$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 10426246
$ws.Range("I132").Value = 4749
$ws.Range("K132").Value = 14247
$ws.Range("M132").Value = -11717

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 18445.455
$ws.Range("I54").Value = 1000
$ws.Range("J54").Value = 24987.5
$ws.Range("K54").Value = 1000
$ws.Range("L54").Value = 24987.5
$ws.Range("M54").Value = -480
$ws.Range("N54").Value = -26027.5

$ws.Range("H107").Value = 525.4
$ws.Range("I107").Value = 409
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 1227
$ws.Range("L107").Value = 2100
$ws.Range("M107").Value = 693
$ws.Range("N107").Value = -5940
